$d = $word.ActiveDocument

$pairs = @(
    @{old="911×2="; new="726×7="},
    @{old="559×8="; new="651×6="},
    @{old="861×8="; new="929×4="},
    @{old="812×9="; new="295×4="},
    @{old="622×7="; new="645×7="},
    @{old="942×5="; new="785×8="},
    @{old="808×2="; new="980×9="},
    @{old="868×8="; new="764×7="},
    @{old="925×4="; new="733×2="},
    @{old="500×2="; new="112×4="},
    @{old="211×5="; new="948×3="},
    @{old="828×8="; new="736×9="},
    @{old="216×4="; new="627×8="},
    @{old="736×2="; new="645×8="},
    @{old="949×3="; new="431×6="},
    @{old="743×8="; new="948×8="},
    @{old="880×8="; new="390×4="},
    @{old="212×3="; new="418×5="},
    @{old="662×7="; new="539×5="},
    @{old="928×5="; new="476×9="},
    @{old="190×6="; new="906×5="},
    @{old="558×9="; new="510×3="},
    @{old="863×4="; new="314×2="},
    @{old="510×2="; new="773×7="},
    @{old="831×4="; new="405×6="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
